$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 11.38889587862907
$ws.Cells.Item(2, 3).Value = 3.716003891676869
$ws.Cells.Item(2, 4).Value = 14.84364476103998
$ws.Cells.Item(2, 5).Value = 16.07547394871409
$ws.Cells.Item(2, 7).Value = 3.710974184978771
$ws.Cells.Item(2, 10).Value = 9.375466445901596
$ws.Cells.Item(2, 11).Value = 10.89353166082071
$ws.Cells.Item(2, 13).Value = 17.22417801385992
$ws.Cells.Item(2, 15).Value = 32.36829310336017

# Row 3
$ws.Cells.Item(3, 2).Value = 11.21319684670871
$ws.Cells.Item(3, 3).Value = 3.591047251481954
$ws.Cells.Item(3, 4).Value = 14.81885661376734
$ws.Cells.Item(3, 5).Value = 16.07365167997235
$ws.Cells.Item(3, 7).Value = 3.713277748980144
$ws.Cells.Item(3, 10).Value = 9.392550704591406
$ws.Cells.Item(3, 11).Value = 10.7896080613981
$ws.Cells.Item(3, 13).Value = 17.19165856557897
$ws.Cells.Item(3, 15).Value = 32.4132195160288

# Row 4
$ws.Cells.Item(4, 2).Value = 11.1067461828113
$ws.Cells.Item(4, 3).Value = 3.512874574355285
$ws.Cells.Item(4, 4).Value = 14.80658339563792
$ws.Cells.Item(4, 5).Value = 16.07540778604723
$ws.Cells.Item(4, 7).Value = 3.71476705038077
$ws.Cells.Item(4, 10).Value = 9.403936452845732
$ws.Cells.Item(4, 11).Value = 10.72769280830856
$ws.Cells.Item(4, 13).Value = 17.17480009864125
$ws.Cells.Item(4, 15).Value = 32.44647363786685

# Row 5
$ws.Cells.Item(5, 2).Value = 11.06378476630651
$ws.Cells.Item(5, 3).Value = 3.480708416683072
$ws.Cells.Item(5, 4).Value = 14.80232664565496
$ws.Cells.Item(5, 5).Value = 16.07684741524767
$ws.Cells.Item(5, 7).Value = 3.715392849671134
$ws.Cells.Item(5, 10).Value = 9.408801848822135
$ws.Cells.Item(5, 11).Value = 10.70296511125361
$ws.Cells.Item(5, 13).Value = 17.16871729102778
$ws.Cells.Item(5, 15).Value = 32.46144858752542

# Row 6
$ws.Cells.Item(6, 2).Value = 11.05667801538048
$ws.Cells.Item(6, 3).Value = 3.475350187172312
$ws.Cells.Item(6, 4).Value = 14.80166488921165
$ws.Cells.Item(6, 5).Value = 16.07713021726285
$ws.Cells.Item(6, 7).Value = 3.715497906306357
$ws.Cells.Item(6, 10).Value = 9.409623379803206
$ws.Cells.Item(6, 11).Value = 10.69889026136702
$ws.Cells.Item(6, 13).Value = 17.16775493586971
$ws.Cells.Item(6, 15).Value = 32.46402108790974

# Row 7
$ws.Cells.Item(7, 2).Value = 11.10616502000933
$ws.Cells.Item(7, 3).Value = 3.512441951306355
$ws.Cells.Item(7, 4).Value = 14.80652296796795
$ws.Cells.Item(7, 5).Value = 16.07542426876808
$ws.Cells.Item(7, 7).Value = 3.71477541353443
$ws.Cells.Item(7, 10).Value = 9.404001155366048
$ws.Cells.Item(7, 11).Value = 10.7273572493128
$ws.Cells.Item(7, 13).Value = 17.17471486983081
$ws.Cells.Item(7, 15).Value = 32.44666983388285

# Row 8
$ws.Cells.Item(8, 2).Value = 11.32805350569116
$ws.Cells.Item(8, 3).Value = 3.673250652417956
$ws.Cells.Item(8, 4).Value = 14.8344885798215
$ws.Cells.Item(8, 5).Value = 16.07424965559794
$ws.Cells.Item(8, 7).Value = 3.711752944771663
$ws.Cells.Item(8, 10).Value = 9.38117133696986
$ws.Cells.Item(8, 11).Value = 10.85732209007979
$ws.Cells.Item(8, 13).Value = 17.21232357142484
$ws.Cells.Item(8, 15).Value = 32.3826060309741

# Row 9
$ws.Cells.Item(9, 2).Value = 11.77193364363671
$ws.Cells.Item(9, 3).Value = 3.974948610965545
$ws.Cells.Item(9, 4).Value = 14.9125246994267
$ws.Cells.Item(9, 5).Value = 16.09469969484258
$ws.Cells.Item(9, 7).Value = 3.706417405147555
$ws.Cells.Item(9, 10).Value = 9.3434969632897
$ws.Cells.Item(9, 11).Value = 11.12598849270929
$ws.Cells.Item(9, 13).Value = 17.31047492002664
$ws.Cells.Item(9, 15).Value = 32.30203562693985

# Row 10
$ws.Cells.Item(10, 2).Value = 12.09996297095863
$ws.Cells.Item(10, 3).Value = 4.185722058014568
$ws.Cells.Item(10, 4).Value = 14.98371486751316
$ws.Cells.Item(10, 5).Value = 16.12350015811042
$ws.Cells.Item(10, 7).Value = 3.702854045269906
$ws.Cells.Item(10, 10).Value = 9.320123727700008
$ws.Cells.Item(10, 11).Value = 11.33010986137211
$ws.Cells.Item(10, 13).Value = 17.39707270339809
$ws.Cells.Item(10, 15).Value = 32.27040010031754

# Row 11
$ws.Cells.Item(11, 2).Value = 12.24890074890502
$ws.Cells.Item(11, 3).Value = 4.278790876109453
$ws.Cells.Item(11, 4).Value = 15.01903784427362
$ws.Cells.Item(11, 5).Value = 16.1395663486052
$ws.Cells.Item(11, 7).Value = 3.701309587504589
$ws.Cells.Item(11, 10).Value = 9.310421616238227
$ws.Cells.Item(11, 11).Value = 11.4240658845577
$ws.Cells.Item(11, 13).Value = 17.43951900390551
$ws.Cells.Item(11, 15).Value = 32.26200540798185

# Row 12
$ws.Cells.Item(12, 2).Value = 12.30520201686902
$ws.Cells.Item(12, 3).Value = 4.313596820219392
$ws.Cells.Item(12, 4).Value = 15.03282921612019
$ws.Cells.Item(12, 5).Value = 16.14607359257806
$ws.Cells.Item(12, 7).Value = 3.700735683117044
$ws.Cells.Item(12, 10).Value = 9.306881159404567
$ws.Cells.Item(12, 11).Value = 11.45977197557943
$ws.Cells.Item(12, 13).Value = 17.45602233183083
$ws.Cells.Item(12, 15).Value = 32.25968920909016

# Row 13
$ws.Cells.Item(13, 2).Value = 12.29308179101684
$ws.Cells.Item(13, 3).Value = 4.306120661744369
$ws.Cells.Item(13, 4).Value = 15.02984064573825
$ws.Cells.Item(13, 5).Value = 16.14465336502556
$ws.Cells.Item(13, 7).Value = 3.700858797672261
$ws.Cells.Item(13, 10).Value = 9.307637726188771
$ws.Cells.Item(13, 11).Value = 11.45207687764618
$ws.Cells.Item(13, 13).Value = 17.45244907858869
$ws.Cells.Item(13, 15).Value = 32.26014967318405

# Row 14
$ws.Cells.Item(14, 2).Value = 12.25353497784775
$ws.Cells.Item(14, 3).Value = 4.281663317789121
$ws.Cells.Item(14, 4).Value = 15.02016418801641
$ws.Cells.Item(14, 5).Value = 16.14009323779398
$ws.Cells.Item(14, 7).Value = 3.70126215295111
$ws.Cells.Item(14, 10).Value = 9.310127666355992
$ws.Cells.Item(14, 11).Value = 11.42700107944392
$ws.Cells.Item(14, 13).Value = 17.44086818217635
$ws.Cells.Item(14, 15).Value = 32.26179756235351

# Row 15
$ws.Cells.Item(15, 2).Value = 12.22929695365952
$ws.Cells.Item(15, 3).Value = 4.266624647841448
$ws.Cells.Item(15, 4).Value = 15.01429094372939
$ws.Cells.Item(15, 5).Value = 16.13735506290011
$ws.Cells.Item(15, 7).Value = 3.701510643804815
$ws.Cells.Item(15, 10).Value = 9.311670206764983
$ws.Cells.Item(15, 11).Value = 11.41165704219113
$ws.Cells.Item(15, 13).Value = 17.43383024640018
$ws.Cells.Item(15, 15).Value = 32.26291929521402

# Row 16
$ws.Cells.Item(16, 2).Value = 12.09021886546288
$ws.Cells.Item(16, 3).Value = 4.179580303778765
$ws.Cells.Item(16, 4).Value = 14.98146495047688
$ws.Cells.Item(16, 5).Value = 16.12250960429372
$ws.Cells.Item(16, 7).Value = 3.70295651436146
$ws.Cells.Item(16, 10).Value = 9.320776483348126
$ws.Cells.Item(16, 11).Value = 11.32398908090351
$ws.Cells.Item(16, 13).Value = 17.39435939556072
$ws.Cells.Item(16, 15).Value = 32.27106942053324

# Row 17
$ws.Cells.Item(17, 2).Value = 12.00478468037548
$ws.Cells.Item(17, 3).Value = 4.125436434038006
$ws.Cells.Item(17, 4).Value = 14.96207487445965
$ws.Cells.Item(17, 5).Value = 16.11415971035294
$ws.Cells.Item(17, 7).Value = 3.7038630702148
$ws.Cells.Item(17, 10).Value = 9.326601010453961
$ws.Cells.Item(17, 11).Value = 11.27046716032087
$ws.Cells.Item(17, 13).Value = 17.37092081859801
$ws.Cells.Item(17, 15).Value = 32.27760550352881

# Row 18
$ws.Cells.Item(18, 2).Value = 11.95562173639777
$ws.Cells.Item(18, 3).Value = 4.094031230397876
$ws.Cells.Item(18, 4).Value = 14.95119923897242
$ws.Cells.Item(18, 5).Value = 16.10963632680219
$ws.Cells.Item(18, 7).Value = 3.704391704147791
$ws.Cells.Item(18, 10).Value = 9.330038723474498
$ws.Cells.Item(18, 11).Value = 11.23978818773919
$ws.Cells.Item(18, 13).Value = 17.35772746811881
$ws.Cells.Item(18, 15).Value = 32.28192930062975

# Row 19
$ws.Cells.Item(19, 2).Value = 11.93897375948886
$ws.Cells.Item(19, 3).Value = 4.083353842766338
$ws.Cells.Item(19, 4).Value = 14.94756472471475
$ws.Cells.Item(19, 5).Value = 16.10815283466248
$ws.Cells.Item(19, 7).Value = 3.70457193006694
$ws.Cells.Item(19, 10).Value = 9.331217727808788
$ws.Cells.Item(19, 11).Value = 11.22941987420043
$ws.Cells.Item(19, 13).Value = 17.35331014228864
$ws.Cells.Item(19, 15).Value = 32.283490184496

# Row 20
$ws.Cells.Item(20, 2).Value = 12.01388214519586
$ws.Cells.Item(20, 3).Value = 4.131227618318841
$ws.Cells.Item(20, 4).Value = 14.96411036060302
$ws.Cells.Item(20, 5).Value = 16.11501968961114
$ws.Cells.Item(20, 7).Value = 3.703765820315063
$ws.Cells.Item(20, 10).Value = 9.325971915454247
$ws.Cells.Item(20, 11).Value = 11.27615397611395
$ws.Cells.Item(20, 13).Value = 17.37338616299705
$ws.Cells.Item(20, 15).Value = 32.27685130733816

# Row 21
$ws.Cells.Item(21, 2).Value = 12.26515394812479
$ws.Cells.Item(21, 3).Value = 4.288859139498731
$ws.Cells.Item(21, 4).Value = 15.02299518626428
$ws.Cells.Item(21, 5).Value = 16.14142119435364
$ws.Cells.Item(21, 7).Value = 3.701143381048374
$ws.Cells.Item(21, 10).Value = 9.30939268906606
$ws.Cells.Item(21, 11).Value = 11.43436325090003
$ws.Cells.Item(21, 13).Value = 17.44425818183507
$ws.Cells.Item(21, 15).Value = 32.26129012286327

# Row 22
$ws.Cells.Item(22, 2).Value = 12.42877694322916
$ws.Cells.Item(22, 3).Value = 4.389319605849402
$ws.Cells.Item(22, 4).Value = 15.06389696215117
$ws.Cells.Item(22, 5).Value = 16.16114193840687
$ws.Cells.Item(22, 7).Value = 3.699493253998253
$ws.Cells.Item(22, 10).Value = 9.299335344251949
$ws.Cells.Item(22, 11).Value = 11.5384879895218
$ws.Cells.Item(22, 13).Value = 17.49307795654217
$ws.Cells.Item(22, 15).Value = 32.25614832211069

# Row 23
$ws.Cells.Item(23, 2).Value = 12.34152134377759
$ws.Cells.Item(23, 3).Value = 4.335946104749288
$ws.Cells.Item(23, 4).Value = 15.04184825714446
$ws.Cells.Item(23, 5).Value = 16.15039206189229
$ws.Cells.Item(23, 7).Value = 3.700368139933135
$ws.Cells.Item(23, 10).Value = 9.304632032093298
$ws.Cells.Item(23, 11).Value = 11.48285850675721
$ws.Cells.Item(23, 13).Value = 17.46679624574451
$ws.Cells.Item(23, 15).Value = 32.25843246878387

# Row 24
$ws.Cells.Item(24, 2).Value = 12.00976931927395
$ws.Cells.Item(24, 3).Value = 4.128610284713764
$ws.Cells.Item(24, 4).Value = 14.96318926960546
$ws.Cells.Item(24, 5).Value = 16.1146300296883
$ws.Cells.Item(24, 7).Value = 3.703809763805444
$ws.Cells.Item(24, 10).Value = 9.326256051693914
$ws.Cells.Item(24, 11).Value = 11.27358267931483
$ws.Cells.Item(24, 13).Value = 17.37227070232522
$ws.Cells.Item(24, 15).Value = 32.27719051598276

# Row 25
$ws.Cells.Item(25, 2).Value = 11.65127216205544
$ws.Cells.Item(25, 3).Value = 3.895074428011268
$ws.Cells.Item(25, 4).Value = 14.88895759186971
$ws.Cells.Item(25, 5).Value = 16.08673998072567
$ws.Cells.Item(25, 7).Value = 3.7077978934371
$ws.Cells.Item(25, 10).Value = 9.352931296220177
$ws.Cells.Item(25, 11).Value = 11.05200118723627
$ws.Cells.Item(25, 13).Value = 17.31047492002664
$ws.Cells.Item(25, 15).Value = 32.31899753492481
